# Update "sodan.xlsx" data: append a new daily row of counts (row 118),
# pushing the trailing footnote row down to row 119, and extend the
# print area / dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row above the current last row (118, the footnote row).
# Excel's default Insert shifts the footnote down to row 119 and copies
# the formatting (number formats / styles) down from row 117 into the
# freshly inserted row 118.
$ws.Rows("118:118").Insert(-4121)

# Fill in the new day's figures in row 118.
$ws.Range("A118").Value2 = 43973
$ws.Range("B118").Value2 = 136
$ws.Range("C118").Value2 = 38601
$ws.Range("D118").Value2 = 37
$ws.Range("E118").Value2 = 7801

# Keep the active selection in sync with the new last row.
$ws.Range("B119").Select()

# Extend the print area to cover the new row.
$pa = $wb.Names.Item(1)
$pa.RefersTo = "=相談件数!`$A`$1:`$E`$119"
